$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GPIO")

# --- Header row (row 1) ---
# Column I header changes from "Code" to "Definitions"
$ws.Range("I1").Value = "Definitions"
# New headers for the two new sections (copy the same bold header format as I1)
$ws.Range("I1").Copy() | Out-Null
$ws.Range("K1").PasteSpecial(-4122) | Out-Null
$ws.Range("M1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("K1").Value = "Global variables"
$ws.Range("M1").Value = "Map"

# --- Data rows: add "static const px_gpio_handle_t ..." (col K) and
#     "#define PX_GPIO_<port><pin> <name>" (col M) formulas ---
$dataRows = @(2,3,4,5,6,7,8,10,11,12,13,14,15,17,18,19,20,21,22,23,24)

foreach ($r in $dataRows) {
    $ws.Range("K$r").Formula = '=CONCATENATE("static const px_gpio_handle_t ",LOWER($A' + $r + '),' + '" = {",$A' + $r + ',"};")'
    $ws.Range("M$r").Formula = '=CONCATENATE("#define PX_GPIO_",$B' + $r + ',$C' + $r + ',"' + ' ' + '",$A' + $r + ')'
}

# --- Column widths for the two new data columns (K and M) ---
$ws.Columns.Item(11).ColumnWidth = 53.833333333333336
$ws.Columns.Item(13).ColumnWidth = 30.666666666666668

Write-Output "done"
